$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54 (shifts old rows 54-82 down to 55-83)
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new record's data.
# (Columns A,B,C,E,F,G,H,I,J,K,L,Q,T are constant across this data block,
#  matching the surrounding rows; D/M/N/O/P/R/S carry the new values.)
$ws.Cells.Item(54, 1).Value = 5
$ws.Cells.Item(54, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(54, 3).Value = "Maule"
$ws.Cells.Item(54, 4).Value = 44489
$ws.Cells.Item(54, 5).Value = 7
$ws.Cells.Item(54, 6).Value = "Fruta"
$ws.Cells.Item(54, 7).Value = 100108
$ws.Cells.Item(54, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(54, 9).Value = 100108002
$ws.Cells.Item(54, 10).Value = "Mango"
$ws.Cells.Item(54, 11).Value = "Sin especificar"
$ws.Cells.Item(54, 12).Value = "Primera"
$ws.Cells.Item(54, 13).Value = 210
$ws.Cells.Item(54, 14).Value = 7000
$ws.Cells.Item(54, 15).Value = 7000
$ws.Cells.Item(54, 16).Value = 7000
$ws.Cells.Item(54, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(54, 18).Value = "Perú"
$ws.Cells.Item(54, 19).Value = 1750
$ws.Cells.Item(54, 20).Value = 4
